$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> column -> new value (from diff). Columns: D = Price, E = Volume(1h).
# ForceText = $true marks Price values that look like plain numbers
# (e.g. "1.00", "510.28") so we can store them as text and keep their
# exact original formatting (trailing zeros, no float drift, no
# scientific notation) instead of Excel auto-converting them to numbers.
$changes = @{
    2 = @{ D = @{ Value = "56.709.90"; ForceText = $false }; E = @{ Value = "  +3.20%  "; ForceText = $false } }
    3 = @{ D = @{ Value = "3.002.51"; ForceText = $false }; E = @{ Value = "  +2.99%  "; ForceText = $false } }
    4 = @{ D = @{ Value = "1.00"; ForceText = $true }; E = @{ Value = "  -0.05%  "; ForceText = $false } }
    5 = @{ D = @{ Value = "510.28"; ForceText = $true }; E = @{ Value = "  +7.40%  "; ForceText = $false } }
    6 = @{ D = @{ Value = "139.77"; ForceText = $true }; E = @{ Value = "  +8.93%  "; ForceText = $false } }
    7 = @{ D = @{ Value = "1.00"; ForceText = $true }; E = @{ Value = "  +0.08%  "; ForceText = $false } }
    8 = @{ E = @{ Value = "  +5.56%  "; ForceText = $false } }
    9 = @{ D = @{ Value = "7.57"; ForceText = $true }; E = @{ Value = "  +12.89%  "; ForceText = $false } }
    10 = @{ D = @{ Value = "0.109"; ForceText = $true }; E = @{ Value = "  +10.95%  "; ForceText = $false } }
    11 = @{ E = @{ Value = "  +5.15%  "; ForceText = $false } }
    12 = @{ E = @{ Value = "  +4.43%  "; ForceText = $false } }
    13 = @{ D = @{ Value = "3.517.03"; ForceText = $false }; E = @{ Value = "  +2.82%  "; ForceText = $false } }
    14 = @{ D = @{ Value = "25.85"; ForceText = $true }; E = @{ Value = "  +9.29%  "; ForceText = $false } }
    15 = @{ E = @{ Value = "  +15.33%  "; ForceText = $false } }
    16 = @{ D = @{ Value = "56.783.52"; ForceText = $false }; E = @{ Value = "  +3.43%  "; ForceText = $false } }
    17 = @{ D = @{ Value = "3.005.87"; ForceText = $false }; E = @{ Value = "  +3.15%  "; ForceText = $false } }
    18 = @{ D = @{ Value = "5.94"; ForceText = $true }; E = @{ Value = "  +9.12%  "; ForceText = $false } }
    19 = @{ D = @{ Value = "12.50"; ForceText = $true }; E = @{ Value = "  +7.87%  "; ForceText = $false } }
    20 = @{ D = @{ Value = "7.86"; ForceText = $true }; E = @{ Value = "  +9.62%  "; ForceText = $false } }
    21 = @{ D = @{ Value = "330.09"; ForceText = $true }; E = @{ Value = "  +8.16%  "; ForceText = $false } }
    22 = @{ E = @{ Value = "  -0.27%  "; ForceText = $false } }
    23 = @{ D = @{ Value = "0.484"; ForceText = $true }; E = @{ Value = "  +8.11%  "; ForceText = $false } }
    24 = @{ D = @{ Value = "62.82"; ForceText = $true }; E = @{ Value = "  +5.43%  "; ForceText = $false } }
    25 = @{ D = @{ Value = "0.172"; ForceText = $true }; E = @{ Value = "  +12.24%  "; ForceText = $false } }
    26 = @{ E = @{ Value = "  +1.00%  "; ForceText = $false } }
    27 = @{ D = @{ Value = "0.0₃0914"; ForceText = $false }; E = @{ Value = "  +11.90%  "; ForceText = $false } }
    28 = @{ D = @{ Value = "6.69"; ForceText = $true }; E = @{ Value = "  +8.09%  "; ForceText = $false } }
    29 = @{ D = @{ Value = "7.15"; ForceText = $true }; E = @{ Value = "  +12.94%  "; ForceText = $false } }
    30 = @{ D = @{ Value = "1.27"; ForceText = $true }; E = @{ Value = "  +12.21%  "; ForceText = $false } }
    31 = @{ E = @{ Value = "  +8.51%  "; ForceText = $false } }
    32 = @{ D = @{ Value = "20.70"; ForceText = $true }; E = @{ Value = "  +9.26%  "; ForceText = $false } }
    33 = @{ D = @{ Value = "156.00"; ForceText = $true }; E = @{ Value = "  +8.02%  "; ForceText = $false } }
    34 = @{ E = @{ Value = "  +7.52%  "; ForceText = $false } }
    35 = @{ D = @{ Value = "5.67"; ForceText = $true }; E = @{ Value = "  +3.71%  "; ForceText = $false } }
    36 = @{ D = @{ Value = "1.27"; ForceText = $true }; E = @{ Value = "  +3.89%  "; ForceText = $false } }
    37 = @{ D = @{ Value = "0.0685"; ForceText = $true }; E = @{ Value = "  +8.97%  "; ForceText = $false } }
    38 = @{ D = @{ Value = "24.36"; ForceText = $true }; E = @{ Value = "  +4.23%  "; ForceText = $false } }
    39 = @{ D = @{ Value = "3.035.85"; ForceText = $false }; E = @{ Value = "  +3.01%  "; ForceText = $false } }
    40 = @{ D = @{ Value = "36.94"; ForceText = $true }; E = @{ Value = "  +3.75%  "; ForceText = $false } }
    41 = @{ E = @{ Value = "  -0.05%  "; ForceText = $false } }
    42 = @{ D = @{ Value = "0.648"; ForceText = $true }; E = @{ Value = "  +5.83%  "; ForceText = $false } }
    43 = @{ D = @{ Value = "2.270.86"; ForceText = $false }; E = @{ Value = "  +10.24%  "; ForceText = $false } }
    44 = @{ D = @{ Value = "1.41"; ForceText = $true }; E = @{ Value = "  +5.42%  "; ForceText = $false } }
    45 = @{ D = @{ Value = "0.998"; ForceText = $true }; E = @{ Value = "  +3.24%  "; ForceText = $false } }
    46 = @{ E = @{ Value = "  +5.46%  "; ForceText = $false } }
    47 = @{ D = @{ Value = "2.00"; ForceText = $true }; E = @{ Value = "  +23.48%  "; ForceText = $false } }
    48 = @{ E = @{ Value = "  +9.51%  "; ForceText = $false } }
    49 = @{ E = @{ Value = "  +7.40%  "; ForceText = $false } }
    50 = @{ D = @{ Value = "19.34"; ForceText = $true }; E = @{ Value = "  +7.02%  "; ForceText = $false } }
    51 = @{ E = @{ Value = "  +9.06%  "; ForceText = $false } }
}

foreach ($row in $changes.Keys) {
    $rowData = $changes[$row]
    foreach ($col in $rowData.Keys) {
        $entry = $rowData[$col]
        $cell = $ws.Range("$col$row")
        if ($entry.ForceText) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $entry.Value
    }
}

Write-Host "Applied cryptos list update"
